$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("MAE") before the existing "Tipo" column, shifting
# "Tipo" (and its data, "multiple") to column E.
$ws.Columns.Item(4).Insert()

# Give the new header cell the same (bold / bordered / centered) formatting
# as the rest of the header row, then set its text.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "MAE"

# MAE values for each disease row.
$ws.Range("D2").Value = 0.5156447814141927
$ws.Range("D3").Value = 0.2551205522039418
$ws.Range("D4").Value = 0.1308862606550939
$ws.Range("D5").Value = 0.2173766023767729

# Updated MSE (B) and R2 (C) values from the retrained model.
$ws.Range("B2").Value = 0.4935191168284981
$ws.Range("C2").Value = 0.9855040507337514

$ws.Range("B3").Value = 0.09660934663994822
$ws.Range("C3").Value = 0.998664832671067

$ws.Range("B4").Value = 0.03319235045495682
$ws.Range("C4").Value = 0.9996578935213501

$ws.Range("B5").Value = 0.09415282550293474
$ws.Range("C5").Value = 0.99943915848242
